$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.725964784622192
$ws.Range("B1").Value = 1.872593283653259
$ws.Range("C1").Value = 2.131636381149292
$ws.Range("D1").Value = 3.399286031723022
$ws.Range("E1").Value = 2.760457515716553
